$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 68, shifting existing rows 68-131 down to 69-132.
$ws.Rows.Item(68).Insert()

# Populate the newly inserted row 68 with the new weekly record.
# Columns that stay identical to the (old) row 68 / new row 69 record:
$ws.Range("A68").Value = 11
$ws.Range("B68").Value = "Vega Monumental Concepción"
$ws.Range("C68").Value = "Bíobío"
$ws.Range("D68").Value = 44729
$ws.Range("E68").Value = 8
$ws.Range("F68").Value = 100112032
$ws.Range("G68").Value = "Zapallo italiano"
$ws.Range("H68").Value = "Sin especificar"
$ws.Range("I68").Value = "Primera"
$ws.Range("J68").Value = 180
$ws.Range("K68").Value = 13000
$ws.Range("L68").Value = 14000
$ws.Range("M68").Value = 13556
$ws.Range("N68").Value = "$/caja 50 unidades"
$ws.Range("O68").Value = "Región de Arica y Parinacota"
$ws.Range("P68").Value = 271
$ws.Range("Q68").Value = 50
$ws.Range("R68").Value = "Hortaliza"
